# LeetCode-Neetcode.xlsx — add three new problem sections (Strings,
# Hashmap & Hashset, Random) below the existing "Two Pointers" section.
#
# Each section is: a bold/italic category header row, a blank spacer row,
# then a hyperlinked problem-URL row, followed by two more blank spacer
# rows before the next header. This mirrors the existing "Two Pointers"
# block (row 3 header, rows 5-6 links).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section: Strings (header row 9, link row 10) -----------------------
$ws.Range("A9").Value = "Strings"
$ws.Range("A3").Copy()
$ws.Range("A9").PasteSpecial(-4122)   # xlPasteFormats: reuse header style

$ws.Hyperlinks.Add($ws.Range("A10"), "https://leetcode.com/problems/length-of-last-word/")

# --- Section: Hashmap & Hashset (header row 13, link row 14) ------------
$ws.Range("A13").Value = "Hashmap & Hashset"
$ws.Range("A3").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("A14"), "https://leetcode.com/problems/contains-duplicate/")

# --- Section: Random (header row 18, link row 19) ------------------------
$ws.Range("A18").Value = "Random"
$ws.Range("A3").Copy()
$ws.Range("A18").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("A19"), "https://leetcode.com/problems/concatenation-of-array/description/")

# Newly-added hyperlink cells should carry the same "Hyperlink" cell style
# already used by the pre-existing links (rows 5 & 6).
$ws.Range("A10").Style = "Hyperlink"
$ws.Range("A14").Style = "Hyperlink"
$ws.Range("A19").Style = "Hyperlink"

# Move the active selection past the new content, matching where the
# author's cursor ended up after typing the new rows.
$ws.Range("A22").Select()

Write-Host "done"
